$wb = $excel.ActiveWorkbook

# Rename the "Include from ..." sheets to the generic "Include #N" names
$wb.Worksheets.Item("Include from MedComCorePracti").Name = "Include #0"
$wb.Worksheets.Item("Include from MedComCorePracti 2").Name = "Include #1"
$wb.Worksheets.Item("Include from NullFlavor").Name = "Include #2"

# Bump the published IG version shown on the Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.8.1"
